# #5: property aircraft done
#
# The "建物" (building) worksheet's property_category column (I) was
# mistakenly left as "land" for all of its data rows; fix it to "building".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("建物")

$ws.Range("I2:I5").Value = "building"
